$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.604.49"
$ws.Range("E2").Value = "  +3.24%  "
$ws.Range("D3").Value = "1.858.17"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "272.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5276"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3375"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06787"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.45%  "
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7915"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07724"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.99%  "
$ws.Range("D13").Value = "1.888.50"
$ws.Range("E13").Value = "  +3.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "89.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.125"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  +1.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007975"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("D20").Value = "26.638.64"
$ws.Range("E20").Value = "  +3.20%  "
$ws.Range("D21").Value = "2.119.66"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.717"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.971"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.095"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.347"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.58%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.653"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.304"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.291"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  +1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04896"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.157"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7257"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.879"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.222"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.315"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01841"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5070"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9382"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "115.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.121"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.995"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4395"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.49%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1321"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.288"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05933"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.469"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
